$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3512
$ws.Range("I64").Value = 3274
$ws.Range("J64").Value = 3750
$ws.Range("K64").Value = 3274
$ws.Range("L64").Value = 3750
$ws.Range("M64").Value = -3026
$ws.Range("N64").Value = -4246
$ws.Range("H67").Value = 3512
$ws.Range("I67").Value = 3274
$ws.Range("J67").Value = 3750
$ws.Range("K67").Value = 3274
$ws.Range("L67").Value = 3750
$ws.Range("M67").Value = -2416
$ws.Range("N67").Value = -5466
$ws.Range("H141").Value = 1447.5758
$ws.Range("I141").Value = 1434.5161
$ws.Range("J141").Value = 1650
$ws.Range("K141").Value = 4303.5483
$ws.Range("L141").Value = 4950
$ws.Range("M141").Value = 876.4516999999996
$ws.Range("N141").Value = -15310

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1744.6111
$ws.Range("I2").Value = 2136.3845
$ws.Range("J2").Value = 726
$ws.Range("K2").Value = 2136.3845
$ws.Range("L2").Value = 726
$ws.Range("M2").Value = -2023.3845
$ws.Range("N2").Value = -952
$ws.Range("H45").Value = 1509
$ws.Range("I45").Value = 1370.1538
$ws.Range("K45").Value = 1370.1538
$ws.Range("M45").Value = -993.1538
$ws.Range("H61").Value = 609433.8
$ws.Range("I61").Value = 772003.4399999999
$ws.Range("K61").Value = 772003.4399999999
$ws.Range("M61").Value = -771791.4399999999
$ws.Range("H74").Value = 191714.02
$ws.Range("I74").Value = 257848.77
$ws.Range("J74").Value = 55963.74
$ws.Range("K74").Value = 257848.77
$ws.Range("L74").Value = 55963.74
$ws.Range("M74").Value = -256974.77
$ws.Range("N74").Value = -57711.74
$ws.Range("H77").Value = 191714.02
$ws.Range("I77").Value = 257848.77
$ws.Range("J77").Value = 55963.74
$ws.Range("K77").Value = 1289243.85
$ws.Range("L77").Value = 279818.7
$ws.Range("M77").Value = -1284875.85
$ws.Range("N77").Value = -288554.7
$ws.Range("H116").Value = 1744.6111
$ws.Range("I116").Value = 2136.3845
$ws.Range("J116").Value = 726
$ws.Range("K116").Value = 2136.3845
$ws.Range("L116").Value = 726
$ws.Range("M116").Value = 157.6154999999999
$ws.Range("N116").Value = -5314
$ws.Range("H132").Value = 2499.2222
$ws.Range("I132").Value = 2188
$ws.Range("J132").Value = 3911.6924
$ws.Range("K132").Value = 6564
$ws.Range("L132").Value = 11735.0772
$ws.Range("M132").Value = -4034
$ws.Range("N132").Value = -16795.0772
$ws.Range("H136").Value = 609433.8
$ws.Range("I136").Value = 772003.4399999999
$ws.Range("K136").Value = 2316010.32
$ws.Range("M136").Value = -2313460.32

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1744.6111
$ws.Range("I3").Value = 2136.3845
$ws.Range("J3").Value = 726
$ws.Range("K3").Value = 2136.3845
$ws.Range("L3").Value = 726
$ws.Range("M3").Value = -2022.3845
$ws.Range("N3").Value = -954
$ws.Range("H64").Value = 540.0909
$ws.Range("J64").Value = 595.55554
$ws.Range("L64").Value = 595.55554
$ws.Range("N64").Value = -1045.55554
$ws.Range("H67").Value = 540.0909
$ws.Range("J67").Value = 595.55554
$ws.Range("L67").Value = 595.55554
$ws.Range("N67").Value = -2155.55554
$ws.Range("H134").Value = 30647.543
$ws.Range("I134").Value = 36691.07
$ws.Range("J134").Value = 6473.4287
$ws.Range("K134").Value = 110073.21
$ws.Range("L134").Value = 19420.2861
$ws.Range("M134").Value = -107538.21
$ws.Range("N134").Value = -24490.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 4375
$ws.Range("I25").Value = 3750
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 3750
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = -3576
$ws.Range("N25").Value = -5348
$ws.Range("H58").Value = 3519.85
$ws.Range("I58").Value = 4006.7026
$ws.Range("J58").Value = 2736.652
$ws.Range("K58").Value = 4006.7026
$ws.Range("L58").Value = 2736.652
$ws.Range("M58").Value = -3803.7026
$ws.Range("N58").Value = -3142.652
$ws.Range("H94").Value = 7376.846
$ws.Range("I94").Value = 1166.6666
$ws.Range("J94").Value = 9239.9
$ws.Range("K94").Value = 1166.6666
$ws.Range("L94").Value = 9239.9
$ws.Range("M94").Value = -715.6666
$ws.Range("N94").Value = -10141.9
$ws.Range("H134").Value = 1490.8085
$ws.Range("I134").Value = 823.8929000000001
$ws.Range("J134").Value = 2473.6316
$ws.Range("K134").Value = 2471.6787
$ws.Range("L134").Value = 7420.8948
$ws.Range("M134").Value = 63.32129999999961
$ws.Range("N134").Value = -12490.8948
$ws.Range("H136").Value = 3519.85
$ws.Range("I136").Value = 4006.7026
$ws.Range("J136").Value = 2736.652
$ws.Range("K136").Value = 12020.1078
$ws.Range("L136").Value = 8209.956
$ws.Range("M136").Value = -9470.1078
$ws.Range("N136").Value = -13309.956
$ws.Range("H141").Value = 25398
$ws.Range("J141").Value = 24424.533
$ws.Range("L141").Value = 24424.533
$ws.Range("N141").Value = -34784.533

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 909.2222
$ws.Range("I113").Value = 1148.5
$ws.Range("J113").Value = 717.8
$ws.Range("K113").Value = 3445.5
$ws.Range("L113").Value = 2153.4
$ws.Range("M113").Value = -1275.5
$ws.Range("N113").Value = -6493.4
$ws.Range("H122").Value = 781.45
$ws.Range("I122").Value = 385.1
$ws.Range("J122").Value = 1177.8
$ws.Range("K122").Value = 3465.9
$ws.Range("L122").Value = 10600.2
$ws.Range("M122").Value = -1015.9
$ws.Range("N122").Value = -15500.2
$ws.Range("H132").Value = 7800
$ws.Range("I132").Value = 6080
$ws.Range("J132").Value = 8755.556
$ws.Range("K132").Value = 54720
$ws.Range("L132").Value = 78800.004
$ws.Range("M132").Value = -52190
$ws.Range("N132").Value = -83860.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 830.4054
$ws.Range("I93").Value = 818.3103599999999
$ws.Range("J93").Value = 874.25
$ws.Range("K93").Value = 818.3103599999999
$ws.Range("L93").Value = 874.25
$ws.Range("M93").Value = 429.6896400000001
$ws.Range("N93").Value = -3370.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 26000
$ws.Range("J48").Value = 26000
$ws.Range("L48").Value = 26000
$ws.Range("N48").Value = -27138
$ws.Range("H122").Value = 2320.8235
$ws.Range("I122").Value = 2278.375
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6835.125
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4385.125
$ws.Range("N122").Value = -13900
